# Update "想去人数" (F column) figures across the three affected sheets
# to match the refreshed data snapshot ("Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1137
$ws1.Range("F6").Value  = 472
$ws1.Range("F7").Value  = 776
$ws1.Range("F8").Value  = 266
$ws1.Range("F11").Value = 427
$ws1.Range("F14").Value = 930
$ws1.Range("F15").Value = 125
$ws1.Range("F16").Value = 2042
$ws1.Range("F17").Value = 524
$ws1.Range("F18").Value = 8537
$ws1.Range("F19").Value = 784
$ws1.Range("F20").Value = 527
$ws1.Range("F25").Value = 142

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5588
$ws3.Range("F4").Value = 407

# --- Sheet "全部类型" (All types, aggregate of the other sheets) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5588
$ws4.Range("F5").Value  = 407
$ws4.Range("F7").Value  = 1137
$ws4.Range("F11").Value = 472
$ws4.Range("F12").Value = 776
$ws4.Range("F14").Value = 266
$ws4.Range("F18").Value = 427
$ws4.Range("F23").Value = 930
$ws4.Range("F24").Value = 125
$ws4.Range("F27").Value = 2042
$ws4.Range("F28").Value = 524
$ws4.Range("F29").Value = 8537
$ws4.Range("F32").Value = 784
$ws4.Range("F33").Value = 527
